$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 325, pushing the existing rows 325..404 down to 326..405
$ws.Rows("325:325").Insert()

# Populate the newly inserted row 325 with the new weekly record
$ws.Range("A325").Value = 3
$ws.Range("B325").Value = "Femacal de La Calera"
$ws.Range("C325").Value = "Coquimbo"
$ws.Range("D325").Value = 44785
$ws.Range("E325").Value = 5
$ws.Range("F325").Value = 100112040
$ws.Range("G325").Value = "Cilantro"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 190
$ws.Range("K325").Value = 4500
$ws.Range("L325").Value = 5000
$ws.Range("M325").Value = 4763
$ws.Range("N325").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O325").Value = "Provincia de Quillota"
$ws.Range("P325").Value = 1588
$ws.Range("Q325").Value = 3
$ws.Range("R325").Value = "Hortaliza"
